$d = $word.ActiveDocument
$shp = $d.InlineShapes.Item(1)
$rng = $shp.Range
Write-Output ("range text len=" + $rng.Text.Length)
$found = $rng.Find.Execute("declercq.denis", $true, $false, $false, $false, $false, $true, 1, $false, "XXXX", 2)
Write-Output ("find result: " + $found)
Write-Output ("range text after: [" + $rng.Text + "]")
